$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Simple in-place edits preceding the insertion block ---
$ws.Range("F1540").Value = "Medium"
$ws.Range("I1546").Value = "N * LOG(N)"
$ws.Range("I1548").Value = "N * LOG(N)"

# --- Insert new rows (ascending order of final row number) ---
$ws.Rows.Item(1549).Insert()
$ws.Range("A1549").Value = 1708
$ws.Range("C1549").Value = 2
$ws.Range("D1549").Value = "C"
$ws.Range("F1549").Value = "Easy"
$ws.Range("G1549").Value = "Array"
$ws.Range("H1549").Value = "Reset head of sub array"
$ws.Range("I1549").Value = "N"
$ws.Range("B1549").Value = "Largest Subarray Length K "

$ws.Rows.Item(1554).Insert()
$ws.Range("A1554").Value = 1714
$ws.Range("C1554").Value = 4
$ws.Range("D1554").Value = "B"
$ws.Range("F1554").Value = "Hard"
$ws.Range("G1554").Value = "Hashtable"
$ws.Range("H1554").Value = "Memorize the query steps"
$ws.Range("I1554").Value = "N * Sqrt(N)"
$ws.Hyperlinks.Add($ws.Range("B1554"), "https://leetcode.com/problems/sum-of-special-evenly-spaced-elements-in-array", "", "", "Sum Of Special Evenly-Spaced Elements In Array ")

$ws.Rows.Item(1563).Insert()
$ws.Range("A1563").Value = 1724
$ws.Range("C1563").Value = 4
$ws.Range("D1563").Value = "C"
$ws.Range("E1563").Value = "***"
$ws.Range("F1563").Value = "Hard"
$ws.Range("G1563").Value = "Graph"
$ws.Range("H1563").Value = "Union Find with weights memorizied "
$ws.Range("I1563").Value = "N * LOG(N)"
$ws.Hyperlinks.Add($ws.Range("B1563"), "https://leetcode.com/problems/checking-existence-of-edge-length-limited-paths-ii", "", "", "Checking Existence of Edge Length Limited Paths II ")

$ws.Rows.Item(1567).Insert()
$ws.Range("A1567").Value = 1728
$ws.Range("C1567").Value = 6
$ws.Range("D1567").Value = "B"
$ws.Range("F1567").Value = "Hard"
$ws.Range("G1567").Value = "DFS"
$ws.Range("H1567").Value = "limit the steps within 70"
$ws.Range("I1567").Value = "E ^ N"
$ws.Hyperlinks.Add($ws.Range("B1567"), "https://leetcode.com/problems/cat-and-mouse-ii", "", "", "Cat and Mouse II")

$ws.Rows.Item(1568).Insert()
$ws.Range("A1568").Value = 1730
$ws.Range("C1568").Value = 3
$ws.Range("D1568").Value = "C"
$ws.Range("E1568").Value = "*"
$ws.Range("F1568").Value = "Medium"
$ws.Range("G1568").Value = "Graph"
$ws.Range("H1568").Value = "BFS to look for all paths"
$ws.Range("I1568").Value = "N"
$ws.Hyperlinks.Add($ws.Range("B1568"), "https://leetcode.com/problems/shortest-path-to-get-food", "", "", "Shortest Path to Get Food ")

$ws.Rows.Item(1569).Insert()
$ws.Range("A1569").Value = 1732
$ws.Range("C1569").Value = 1
$ws.Range("D1569").Value = "E"
$ws.Range("F1569").Value = "Easy"
$ws.Range("G1569").Value = "Array"
$ws.Range("H1569").Value = "Calculate highest altitude on the way"
$ws.Range("I1569").Value = "N"
$ws.Hyperlinks.Add($ws.Range("B1569"), "https://leetcode.com/problems/find-the-highest-altitude", "", "", "Find the Highest Altitude")

$ws.Rows.Item(1570).Insert()
$ws.Range("A1570").Value = 1736
$ws.Range("C1570").Value = 2
$ws.Range("D1570").Value = "B"
$ws.Range("F1570").Value = "Easy"
$ws.Range("G1570").Value = "String"
$ws.Range("H1570").Value = "Branch Logic"
$ws.Range("I1570").Value = "One"
$ws.Hyperlinks.Add($ws.Range("B1570"), "https://leetcode.com/problems/latest-time-by-replacing-hidden-digits", "", "", "Latest Time by Replacing Hidden Digits ")

$ws.Rows.Item(1571).Insert()
$ws.Range("A1571").Value = 1740
$ws.Range("C1571").Value = 3
$ws.Range("D1571").Value = "C"
$ws.Range("E1571").Value = "**"
$ws.Range("F1571").Value = "Medium"
$ws.Range("G1571").Value = "Tree"
$ws.Range("H1571").Value = "Calculate depth for nodes and calculate result by post-order"
$ws.Range("I1571").Value = "N"
$ws.Hyperlinks.Add($ws.Range("B1571"), "https://leetcode.com/problems/find-distance-in-a-binary-tree", "", "", "Find Distance in a Binary Tree")

$ws.Rows.Item(1572).Insert()
$ws.Range("A1572").Value = 1742
$ws.Range("C1572").Value = 1
$ws.Range("D1572").Value = "E"
$ws.Range("F1572").Value = "Easy"
$ws.Range("G1572").Value = "Math"
$ws.Range("H1572").Value = "Calculate digits"
$ws.Range("I1572").Value = "N"
$ws.Hyperlinks.Add($ws.Range("B1572"), "https://leetcode.com/problems/maximum-number-of-balls-in-a-box", "", "", "Maximum Number of Balls in a Box")

$ws.Rows.Item(1573).Insert()
$ws.Range("A1573").Value = 1746
$ws.Range("C1573").Value = 3
$ws.Range("D1573").Value = "C"
$ws.Range("E1573").Value = "*"
$ws.Range("F1573").Value = "Medium"
$ws.Range("G1573").Value = "DP"
$ws.Range("H1573").Value = "Calculate presum and presum with one operation"
$ws.Range("I1573").Value = "N"
$ws.Hyperlinks.Add($ws.Range("B1573"), "https://leetcode.com/problems/maximum-subarray-sum-after-one-operation", "", "", "Maximum Subarray Sum After One Operation")

$ws.Rows.Item(1574).Insert()
$ws.Range("A1574").Value = 1748
$ws.Range("C1574").Value = 1
$ws.Range("D1574").Value = "E"
$ws.Range("F1574").Value = "Easy"
$ws.Range("G1574").Value = "Hashtable"
$ws.Range("H1574").Value = "Count unique number in array"
$ws.Range("I1574").Value = "N"
$ws.Hyperlinks.Add($ws.Range("B1574"), "https://leetcode.com/problems/sum-of-unique-elements", "", "", "Sum of Unique Elements")

$ws.Rows.Item(1575).Insert()
$ws.Range("A1575").Value = 1752
$ws.Range("C1575").Value = 2
$ws.Range("D1575").Value = "C"
$ws.Range("E1575").Value = "*"
$ws.Range("F1575").Value = "Easy"
$ws.Range("G1575").Value = "Sort"
$ws.Range("H1575").Value = "Count dip and if dip is 1 compare first and last element"
$ws.Range("I1575").Value = "N"
$ws.Hyperlinks.Add($ws.Range("B1575"), "https://leetcode.com/problems/check-if-array-is-sorted-and-rotated", "", "", "Check if Array Is Sorted and Rotated")

$ws.Rows.Item(1576).Insert()
$ws.Range("A1576").Value = 1756
$ws.Range("C1576").Value = 3
$ws.Range("D1576").Value = "C"
$ws.Range("F1576").Value = "Medium"
$ws.Range("G1576").Value = "Design"
$ws.Range("H1576").Value = "Simply move elements in array, a more advanced version will be using BIT"
$ws.Range("I1576").Value = "N * N -> N Log(N)"
$ws.Hyperlinks.Add($ws.Range("B1576"), "https://leetcode.com/problems/design-most-recently-used-queue", "", "", "Design Most Recently Used Queue")
